$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells P1, Q1 - copy formatting (style) from O1, then set values
$ws.Range("O1").Copy($ws.Range("P1:Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For each data row (2-25), swap I<->K and M<->O values (1 <-> 2), and add P/Q = 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I
    $ws.Cells.Item($r, 11).Value = 1  # K
    $ws.Cells.Item($r, 13).Value = 2  # M
    $ws.Cells.Item($r, 15).Value = 1  # O
    $ws.Cells.Item($r, 16).Value = 2  # P
    $ws.Cells.Item($r, 17).Value = 2  # Q
}
